# Applies the crypto price/volume refresh described in the commit diff.
# Only the cells that actually changed between before/after are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Forces the cell to keep/receive a *text* value (not a number/date),
    # matching the inlineStr cells in the original workbook, while
    # restoring the original cell style so no stray styles are introduced.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2: Bitcoin
Set-TextValue $ws.Range("D2") "66.349.36"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3: Ethereum
Set-TextValue $ws.Range("D3") "3.328.53"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4: TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5: Solana
Set-TextValue $ws.Range("D5") "190.66"
$ws.Range("E5").Value = "  +4.20%  "

# Row 6: BNB
Set-TextValue $ws.Range("D6") "563.50"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7: USDC
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8: XRP
Set-TextValue $ws.Range("D8") "0.589"
$ws.Range("E8").Value = "  -1.86%  "

# Row 9: LidoStakedEther
Set-TextValue $ws.Range("D9") "3.324.03"
$ws.Range("E9").Value = "  -0.94%  "

# Row 10: Dogecoin
Set-TextValue $ws.Range("D10") "0.186"
$ws.Range("E10").Value = "  -0.59%  "

# Row 11: Cardano
Set-TextValue $ws.Range("D11") "0.589"
$ws.Range("E11").Value = "  -0.77%  "

# Row 12: Avalanche
Set-TextValue $ws.Range("D12") "48.02"
$ws.Range("E12").Value = "  +0.27%  "

# Row 13: ShibaInu
Set-TextValue $ws.Range("D13") "0.0000273"
$ws.Range("E13").Value = "  +1.22%  "

# Row 14: Polkadot
Set-TextValue $ws.Range("D14") "8.71"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.853.82"
$ws.Range("E15").Value = "  -1.22%  "

# Row 16: BitcoinCash
Set-TextValue $ws.Range("D16") "609.83"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17: Chainlink
Set-TextValue $ws.Range("D17") "18.19"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18: WrappedBTC
Set-TextValue $ws.Range("D18") "66.316.51"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19: TRON
$ws.Range("E19").Value = "  +0.43%  "

# Row 20: WrappedEther
Set-TextValue $ws.Range("D20") "3.307.99"
$ws.Range("E20").Value = "  -1.60%  "

# Row 21: Uniswap
Set-TextValue $ws.Range("D21") "11.16"
$ws.Range("E21").Value = "  -2.93%  "

# Row 22: Polygon
Set-TextValue $ws.Range("D22") "0.917"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D23") "18.52"
$ws.Range("E23").Value = "  +9.83%  "

# Row 24: Toncoin
Set-TextValue $ws.Range("D24") "5.14"
$ws.Range("E24").Value = "  -0.23%  "

# Row 25: Litecoin
Set-TextValue $ws.Range("D25") "101.41"
$ws.Range("E25").Value = "  +1.19%  "

# Row 26: PancakeSwap
Set-TextValue $ws.Range("D26") "4.01"
$ws.Range("E26").Value = "  -1.91%  "

# Row 27: ImmutableX
Set-TextValue $ws.Range("D27") "2.76"
$ws.Range("E27").Value = "  +1.89%  "

# Row 28: RenderToken
Set-TextValue $ws.Range("D28") "9.80"
$ws.Range("E28").Value = "  +4.26%  "

# Row 29: Filecoin
Set-TextValue $ws.Range("D29") "8.71"
$ws.Range("E29").Value = "  -1.01%  "

# Row 30: EthereumClassic
Set-TextValue $ws.Range("D30") "30.49"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31: NEARProtocol
Set-TextValue $ws.Range("D31") "6.80"
$ws.Range("E31").Value = "  +7.69%  "

# Row 32: dogwifhat
Set-TextValue $ws.Range("D32") "4.13"
$ws.Range("E32").Value = "  +7.22%  "

# Row 33: Cosmos
Set-TextValue $ws.Range("D33") "11.16"
$ws.Range("E33").Value = "  +0.25%  "

# Row 34: Bittensor
Set-TextValue $ws.Range("D34") "562.26"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35: Hedera
$ws.Range("E35").Value = "  +0.53%  "

# Row 36: Maker
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D36") "3.729.19"
$ws.Range("E36").Value = "  -3.16%  "

# Row 37: OKB
Set-TextValue $ws.Range("D37") "57.36"
$ws.Range("E37").Value = "  -1.51%  "

# Row 38: Dai
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D38") "1.00"
$ws.Range("E38").Value = "  +0.21%  "

# Row 39: PEPE
Set-TextValue $ws.Range("D39") "0.0₃0733"
$ws.Range("E39").Value = "  +1.67%  "

# Row 40: InjectiveProtocol
Set-TextValue $ws.Range("D40") "34.27"
$ws.Range("E40").Value = "  +6.16%  "

# Row 41: Stacks
Set-TextValue $ws.Range("D41") "3.33"
$ws.Range("E41").Value = "  -2.54%  "

# Row 42: CoreDAO
$ws.Range("B42").Value = "CoreDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue $ws.Range("D42") "3.47"
$ws.Range("E42").Value = "  +1.45%  "

# Row 43: Kaspa
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D43") "0.130"
$ws.Range("E43").Value = "  +2.08%  "

# Row 44: Fetch.AI
Set-TextValue $ws.Range("D44") "2.72"
$ws.Range("E44").Value = "  +2.18%  "

# Row 45: TheGraph
Set-TextValue $ws.Range("D45") "0.344"
$ws.Range("E45").Value = "  -1.15%  "

# Row 46: VeChain
Set-TextValue $ws.Range("D46") "0.0426"
$ws.Range("E46").Value = "  +2.25%  "

# Row 47: ApeXProtocol
Set-TextValue $ws.Range("D47") "3.24"
$ws.Range("E47").Value = "  +2.95%  "

# Row 48: Stellar
$ws.Range("E48").Value = "  -0.37%  "

# Row 49: ThetaToken
Set-TextValue $ws.Range("D49") "2.61"
$ws.Range("E49").Value = "  -1.96%  "

# Row 50: FirstDigitalUSD
Set-TextValue $ws.Range("D50") "0.998"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51: Mantle
Set-TextValue $ws.Range("D51") "1.29"
$ws.Range("E51").Value = "  +2.42%  "

